$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Formats already used in the sheet (column A = date format, column F = currency format).
# Set as literal format codes (rather than round-tripped via .NumberFormat) so the
# non-ASCII currency glyph survives intact and reuses the existing style entries.
$dateFmt = "[$-409]d/mmm/yyyy;@"
$amtFmt = '"₹"#,##0;"₹"\-#,##0'

# Row 197: KA 03 MW 3617 / XCENT / PMS WW / WORK IN PROGRESS
$ws.Cells.Item(197, 1).Value = 44774
$ws.Cells.Item(197, 1).NumberFormat = $dateFmt
$ws.Cells.Item(197, 2).Value = "KA 03 MW 3617"
$ws.Cells.Item(197, 3).Value = "XCENT"
$ws.Cells.Item(197, 4).Value = "PMS                                      WW"
$ws.Cells.Item(197, 5).Value = "WORK IN PROGRESS"

# Row 198: KA 51 MB 4552 / POLO / SUSPENSION / WORK IN PROGRESS
$ws.Cells.Item(198, 1).Value = 44774
$ws.Cells.Item(198, 1).NumberFormat = $dateFmt
$ws.Cells.Item(198, 2).Value = "KA 51 MB 4552"
$ws.Cells.Item(198, 3).Value = "POLO"
$ws.Cells.Item(198, 4).Value = "SUSPENSION"
$ws.Cells.Item(198, 5).Value = "WORK IN PROGRESS"

# Row 199: KA 53 MD 9553 / H CITY / BRAKE DISC & PAD CHANGE / WORK DONE DELIVERED / 2726 / P PAY
$ws.Cells.Item(199, 1).Value = 44774
$ws.Cells.Item(199, 1).NumberFormat = $dateFmt
$ws.Cells.Item(199, 2).Value = "KA 53 MD 9553"
$ws.Cells.Item(199, 3).Value = "H CITY"
$ws.Cells.Item(199, 4).Value = "BRAKE DISC & PAD CHANGE"
$ws.Cells.Item(199, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(199, 6).Value = 2726
$ws.Cells.Item(199, 6).NumberFormat = $amtFmt
$ws.Cells.Item(199, 7).Value = "P PAY"

# Row 200: KL 01 CF 1995 / TIAGO / AC REFLLING / WORK DONE DELIVERED / 2360 / CARD
$ws.Cells.Item(200, 1).Value = 44774
$ws.Cells.Item(200, 1).NumberFormat = $dateFmt
$ws.Cells.Item(200, 2).Value = "KL 01 CF 1995"
$ws.Cells.Item(200, 3).Value = "TIAGO"
$ws.Cells.Item(200, 4).Value = "AC REFLLING "
$ws.Cells.Item(200, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(200, 6).Value = 2360
$ws.Cells.Item(200, 6).NumberFormat = $amtFmt
$ws.Cells.Item(200, 7).Value = "CARD"

# Row 201: KA 03 MT 2522 / ETIOS / BODY SHOP / WORK IN PROGRESS
$ws.Cells.Item(201, 1).Value = 44774
$ws.Cells.Item(201, 1).NumberFormat = $dateFmt
$ws.Cells.Item(201, 2).Value = "KA 03 MT 2522"
$ws.Cells.Item(201, 3).Value = "ETIOS"
$ws.Cells.Item(201, 4).Value = "BODY SHOP"
$ws.Cells.Item(201, 5).Value = "WORK IN PROGRESS"

# Row 202: KA 53 MB 1800 / SCALA / AC REFLLING / WORK DONE DELIVERED / 2000 / P PAY
$ws.Cells.Item(202, 1).Value = 44774
$ws.Cells.Item(202, 1).NumberFormat = $dateFmt
$ws.Cells.Item(202, 2).Value = "KA 53 MB 1800"
$ws.Cells.Item(202, 3).Value = "SCALA"
$ws.Cells.Item(202, 4).Value = "AC REFLLING "
$ws.Cells.Item(202, 5).Value = "WORK DONE DELIVERED"
$ws.Cells.Item(202, 6).Value = 2000
$ws.Cells.Item(202, 6).NumberFormat = $amtFmt
$ws.Cells.Item(202, 7).Value = "P PAY"

# Match the saved view state: scrolled so row 184 is at top, active cell D203
$ws.Application.ActiveWindow.ScrollRow = 184
$ws.Range("D203").Select() | Out-Null
